# Multiple fixes. imported labelmap classes work now
#
# Target sheet "plans" (8th sheet) gets:
#  - row 2: new C2 = 0 (number), new G2 = FALSE (boolean, custom TRUE/FALSE numberformat)
#  - row 3: G3 string "0" -> boolean FALSE ; N3 "TSL.lungs,TSL.lungs" -> "TSL.label_localiser,TSL.lungs"
#  - row 4 (new): plan3 definition mirroring plan2's structure, with a new imported_folder path
#    and the new remapping_imported label
#  - column E and N get custom widths
#  - all sheets: zoom 80%
#  - "plans" sheet selection -> G4

$wb = $excel.ActiveWorkbook

$boolFmt = '"TRUE";"TRUE";"FALSE"'

$ws = $wb.Worksheets.Item("plans")

# ---- Row 2 additions -------------------------------------------------
$ws.Range("C2").Value = 0

$ws.Range("G2").NumberFormat = $boolFmt
$ws.Range("G2").Value = $false

# ---- Row 3 edits -------------------------------------------------------
# G3: was the text "0", becomes boolean FALSE with the custom TRUE/FALSE format
$ws.Range("G3").NumberFormat = $boolFmt
$ws.Range("G3").Value = $false

# N3: remapping_imported label text changes (write this BEFORE the row4 new
# strings so the shared-string table grows in the same order as the source
# commit: TSL.label_localiser,... gets appended first).
$ws.Range("N3").Value = "TSL.label_localiser,TSL.lungs"

# ---- Row 4: brand new "plan3" row --------------------------------------
$ws.Range("A4").Value = "plan3"
$ws.Range("B4").Value = "lidc"

# C4 / I4 / J4 / K4 / P4 hold digit-only text (matching the sibling plan2
# row, which stores them as text, not numbers) -- force text typing with a
# throw-away "@" format so Value isn't auto-coerced to a number, then put
# the format back the way it was.
foreach ($ref in @("C4", "I4", "J4", "K4", "P4")) {
    $ws.Range($ref).NumberFormat = "@"
}
$ws.Range("C4").Value = "0"
$ws.Range("E4").Value = "/s/fran_storage/predictions/totalseg/LITS-1271"
$ws.Range("H4").Value = "lbd"
$ws.Range("I4").Value = "128"
$ws.Range("J4").Value = "96"
$ws.Range("K4").Value = "0.25"
$ws.Range("P4").Value = "2"
foreach ($ref in @("C4", "I4", "J4", "K4", "P4")) {
    $ws.Range($ref).NumberFormat = "General"
}

# L4 / M4 stay empty (like L3/M3) but the cell still needs to exist.
$ws.Range("L4").NumberFormat = "General"
$ws.Range("M4").NumberFormat = "General"

$ws.Range("N4").Value = "TSL.label_localiser,TSL.lungs"
$ws.Range("Q4").Value = "0.8,.8,1.5"
$ws.Range("R4").Value = "manual_value"

$ws.Range("G4").NumberFormat = $boolFmt
$ws.Range("G4").Value = $false

# ---- Column widths ------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 55.5    # column E -> imported_folder
$ws.Columns.Item(14).ColumnWidth = 23.75  # column N -> remapping_imported

# ---- Zoom every sheet to 80% --------------------------------------------
foreach ($sh in $wb.Worksheets) {
    $sh.Activate()
    $excel.ActiveWindow.Zoom = 80
}

# ---- Selection / active sheet -------------------------------------------
$ws.Activate()
$ws.Range("G4").Select()
